$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the question block (rows 2:4) twice, into rows 5:7 and 8:10,
# preserving cell values and formatting (styles), the way Excel's Copy does.
$src = $ws.Range("A2:J4")
$src.Copy($ws.Range("A5:J7"))
$src.Copy($ws.Range("A8:J10"))

# Match the custom row height used by the question rows.
for ($r = 5; $r -le 10; $r++) {
    $ws.Rows.Item($r).RowHeight = 49.2
}

# Extend the "0/1" list validation that was only covering rows 2-4 on
# columns D,F,H,J so that it also covers the newly added rows 5-10.
$ws.Range("D2:D4").Validation.Delete()
$ws.Range("F2:F4").Validation.Delete()
$ws.Range("H2:H4").Validation.Delete()
$ws.Range("J2:J4").Validation.Delete()

$ws.Range("D2:D10").Validation.Add(3, 1, 1, "0,1")
$ws.Range("D2:D10").Validation.IgnoreBlank = $true
$ws.Range("D2:D10").Validation.InCellDropdown = $true
$ws.Range("D2:D10").Validation.ShowInput = $true
$ws.Range("D2:D10").Validation.ShowError = $true

$ws.Range("F2:F10").Validation.Add(3, 1, 1, "0,1")
$ws.Range("F2:F10").Validation.IgnoreBlank = $true
$ws.Range("F2:F10").Validation.InCellDropdown = $true
$ws.Range("F2:F10").Validation.ShowInput = $true
$ws.Range("F2:F10").Validation.ShowError = $true

$ws.Range("H2:H10").Validation.Add(3, 1, 1, "0,1")
$ws.Range("H2:H10").Validation.IgnoreBlank = $true
$ws.Range("H2:H10").Validation.InCellDropdown = $true
$ws.Range("H2:H10").Validation.ShowInput = $true
$ws.Range("H2:H10").Validation.ShowError = $true

$ws.Range("J2:J10").Validation.Add(3, 1, 1, "0,1")
$ws.Range("J2:J10").Validation.IgnoreBlank = $true
$ws.Range("J2:J10").Validation.InCellDropdown = $true
$ws.Range("J2:J10").Validation.ShowInput = $true
$ws.Range("J2:J10").Validation.ShowError = $true

# Update the view: scroll so row 7 is the top row, and select A8:XFD10
# (the newly appended block of rows) as the active selection.
$ws.Range("A8:XFD10").Select()
$excel.ActiveWindow.ScrollRow = 7
